$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired data (player, position, team) in row order after the edit.
$data = @(
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Isaiah Collier", "PG,SG", "Utah Jazz"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Luka Doncic", "PG,SG", "Los Angeles Lakers"),
    @("De'Aaron Fox", "PG", "San Antonio Spurs"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Royce O'Neale", "SF,PF", "Phoenix Suns"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("P.J. Washington", "SF,PF", "Dallas Mavericks")
)

# Clear out the old data rows (rows 2 through 18) before writing the new, shorter table.
$ws.Range("A2:C18").ClearContents()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
